# JTown special collection picture addition first batch
#
# 1. The <freezingUrl> paragraph gets its URL turned into a real
#    hyperlink (the selection/auto-link accidentally swallowed the
#    following "</freezingUrl" text too, exactly like the source diff).
# 2. Two blank paragraphs and a new paragraph containing an imgur
#    picture link are added right after that paragraph.

$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Step 1: turn part of the <freezingUrl> paragraph into a hyperlink
# ---------------------------------------------------------------
$freezingPara = $d.Paragraphs(3)
$paraStart = $freezingPara.Range.Start

# "<freezingUrl>" is the first 13 characters of the paragraph; the
# hyperlink covers the next 58 characters (the URL plus the stray
# "</freezingUrl" text), leaving the closing ">" as plain text.
$linkStart = $paraStart + 13
$linkEnd = $linkStart + 58
$linkRange = $d.Range($linkStart, $linkEnd)

$d.Hyperlinks.Add($linkRange, "https://nchfp.uga.edu/how/freeze/avocado.html")

# Re-apply the named "Hyperlink" character style onto the freshly
# created hyperlink range so Word actually persists the style
# definition into styles.xml (matching what real Word does the first
# time a Hyperlink style gets used in a saved document).
$hlRange = $d.Hyperlinks(1).Range
$hlRange.Style = "Hyperlink"

$hlStyle = $d.Styles("Hyperlink")
$hlStyle.Priority = 99
$hlStyle.UnhideWhenUsed = $true
$hlStyle.QuickStyle = $false

# ---------------------------------------------------------------
# Step 2: two blank paragraphs + a new paragraph with the picture URL
# ---------------------------------------------------------------
$freezingPara = $d.Paragraphs(3)
$freezingPara.Range.InsertParagraphAfter()
$freezingPara.Range.InsertParagraphAfter()
$freezingPara.Range.InsertParagraphAfter()

$picturePara = $d.Paragraphs(6)
$picturePara.Range.InsertAfter("https://i.imgur.com/K4F1Hl8.jpg")

# ---------------------------------------------------------------
# Also materialize the companion "Unresolved Mention" character style
# that Word bundles alongside the Hyperlink style the first time a
# hyperlink-ish style gets unhidden in a saved document (unused by any
# run, but present in styles.xml).
# ---------------------------------------------------------------
$umStyle = $d.Styles.Add("UnresolvedMention", 2)
$umStyle.NameLocal = "Unresolved Mention"
$umStyle.BaseStyle = "DefaultParagraphFont"
$umStyle.Priority = 99
$umStyle.UnhideWhenUsed = $true
# Word COM colors are BGR-ordered (COLORREF); 0x5C5E60 round-trips to
# the target RGB hex "605E5C".
$umStyle.Font.Color = 6053472

Write-Host "Done."
